$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Cells.Item(58, 1).Value = "'2026-01-28"
$ws.Cells.Item(58, 1).ClearFormats()
$ws.Cells.Item(58, 2).Value = "18:31:55"
$ws.Cells.Item(58, 3).Value = "18:00"
$ws.Cells.Item(58, 4).Value = "Bedroom"
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(58, 6).Value = "Occupied"
$ws.Cells.Item(59, 1).Value = "'2026-01-28"
$ws.Cells.Item(59, 1).ClearFormats()
$ws.Cells.Item(59, 2).Value = "18:31:58"
$ws.Cells.Item(59, 3).Value = "18:00"
$ws.Cells.Item(59, 4).Value = "Bedroom"
$ws.Cells.Item(59, 5).Value = 68
$ws.Cells.Item(59, 6).Value = "Occupied"
$ws.Cells.Item(60, 1).Value = "'2026-01-28"
$ws.Cells.Item(60, 1).ClearFormats()
$ws.Cells.Item(60, 2).Value = "18:32:01"
$ws.Cells.Item(60, 3).Value = "18:00"
$ws.Cells.Item(60, 4).Value = "Bedroom"
$ws.Cells.Item(60, 5).Value = 51
$ws.Cells.Item(60, 6).Value = "Occupied"
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 1).ClearFormats()
$ws.Cells.Item(61, 2).Value = "18:32:05"
$ws.Cells.Item(61, 3).Value = "18:00"
$ws.Cells.Item(61, 4).Value = "Bedroom"
$ws.Cells.Item(61, 5).Value = 59
$ws.Cells.Item(61, 6).Value = "Occupied"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 1).ClearFormats()
$ws.Cells.Item(62, 2).Value = "18:32:09"
$ws.Cells.Item(62, 3).Value = "18:00"
$ws.Cells.Item(62, 4).Value = "Bedroom"
$ws.Cells.Item(62, 5).Value = 53
$ws.Cells.Item(62, 6).Value = "Occupied"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 1).ClearFormats()
$ws.Cells.Item(63, 2).Value = "18:32:12"
$ws.Cells.Item(63, 3).Value = "18:00"
$ws.Cells.Item(63, 4).Value = "Bedroom"
$ws.Cells.Item(63, 5).Value = 65
$ws.Cells.Item(63, 6).Value = "Occupied"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 1).ClearFormats()
$ws.Cells.Item(64, 2).Value = "18:32:15"
$ws.Cells.Item(64, 3).Value = "18:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = 54
$ws.Cells.Item(64, 6).Value = "Occupied"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 1).ClearFormats()
$ws.Cells.Item(65, 2).Value = "18:32:18"
$ws.Cells.Item(65, 3).Value = "18:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = 50
$ws.Cells.Item(65, 6).Value = "Occupied"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 1).ClearFormats()
$ws.Cells.Item(66, 2).Value = "18:32:21"
$ws.Cells.Item(66, 3).Value = "18:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = 49
$ws.Cells.Item(66, 6).Value = "Occupied"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 1).ClearFormats()
$ws.Cells.Item(67, 2).Value = "18:32:24"
$ws.Cells.Item(67, 3).Value = "18:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = 57
$ws.Cells.Item(67, 6).Value = "Occupied"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 1).ClearFormats()
$ws.Cells.Item(68, 2).Value = "18:32:27"
$ws.Cells.Item(68, 3).Value = "18:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = 50
$ws.Cells.Item(68, 6).Value = "Occupied"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 1).ClearFormats()
$ws.Cells.Item(69, 2).Value = "18:32:30"
$ws.Cells.Item(69, 3).Value = "18:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = 49
$ws.Cells.Item(69, 6).Value = "Occupied"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 1).ClearFormats()
$ws.Cells.Item(70, 2).Value = "18:32:34"
$ws.Cells.Item(70, 3).Value = "18:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = 50
$ws.Cells.Item(70, 6).Value = "Occupied"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 1).ClearFormats()
$ws.Cells.Item(71, 2).Value = "18:32:37"
$ws.Cells.Item(71, 3).Value = "18:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = 49
$ws.Cells.Item(71, 6).Value = "Occupied"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(72, 2).Value = "18:32:40"
$ws.Cells.Item(72, 3).Value = "18:00"
$ws.Cells.Item(72, 4).Value = "Bedroom"
$ws.Cells.Item(72, 5).Value = 50
$ws.Cells.Item(72, 6).Value = "Occupied"
$ws.Cells.Item(73, 1).Value = "'2026-01-28"
$ws.Cells.Item(73, 1).ClearFormats()
$ws.Cells.Item(73, 2).Value = "18:32:43"
$ws.Cells.Item(73, 3).Value = "18:00"
$ws.Cells.Item(73, 4).Value = "Bedroom"
$ws.Cells.Item(73, 5).Value = 49
$ws.Cells.Item(73, 6).Value = "Occupied"

$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Cells.Item(58, 1).Value = "'2026-01-28"
$ws.Cells.Item(58, 1).ClearFormats()
$ws.Cells.Item(58, 2).Value = "18:31:56"
$ws.Cells.Item(58, 3).Value = "18:00"
$ws.Cells.Item(58, 4).Value = "Bedroom"
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(58, 6).Value = "Occupied"
$ws.Cells.Item(59, 1).Value = "'2026-01-28"
$ws.Cells.Item(59, 1).ClearFormats()
$ws.Cells.Item(59, 2).Value = "18:31:59"
$ws.Cells.Item(59, 3).Value = "18:00"
$ws.Cells.Item(59, 4).Value = "Bedroom"
$ws.Cells.Item(59, 5).Value = 20
$ws.Cells.Item(59, 6).Value = "Occupied"
$ws.Cells.Item(60, 1).Value = "'2026-01-28"
$ws.Cells.Item(60, 1).ClearFormats()
$ws.Cells.Item(60, 2).Value = "18:32:02"
$ws.Cells.Item(60, 3).Value = "18:00"
$ws.Cells.Item(60, 4).Value = "Bedroom"
$ws.Cells.Item(60, 5).Value = 3
$ws.Cells.Item(60, 6).Value = "Occupied"
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 1).ClearFormats()
$ws.Cells.Item(61, 2).Value = "18:32:06"
$ws.Cells.Item(61, 3).Value = "18:00"
$ws.Cells.Item(61, 4).Value = "Bedroom"
$ws.Cells.Item(61, 5).Value = 11
$ws.Cells.Item(61, 6).Value = "Occupied"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 1).ClearFormats()
$ws.Cells.Item(62, 2).Value = "18:32:10"
$ws.Cells.Item(62, 3).Value = "18:00"
$ws.Cells.Item(62, 4).Value = "Bedroom"
$ws.Cells.Item(62, 5).Value = 5
$ws.Cells.Item(62, 6).Value = "Occupied"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 1).ClearFormats()
$ws.Cells.Item(63, 2).Value = "18:32:13"
$ws.Cells.Item(63, 3).Value = "18:00"
$ws.Cells.Item(63, 4).Value = "Bedroom"
$ws.Cells.Item(63, 5).Value = 17
$ws.Cells.Item(63, 6).Value = "Occupied"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 1).ClearFormats()
$ws.Cells.Item(64, 2).Value = "18:32:16"
$ws.Cells.Item(64, 3).Value = "18:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = 6
$ws.Cells.Item(64, 6).Value = "Occupied"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 1).ClearFormats()
$ws.Cells.Item(65, 2).Value = "18:32:19"
$ws.Cells.Item(65, 3).Value = "18:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = 2
$ws.Cells.Item(65, 6).Value = "Occupied"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 1).ClearFormats()
$ws.Cells.Item(66, 2).Value = "18:32:22"
$ws.Cells.Item(66, 3).Value = "18:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = 1
$ws.Cells.Item(66, 6).Value = "Occupied"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 1).ClearFormats()
$ws.Cells.Item(67, 2).Value = "18:32:25"
$ws.Cells.Item(67, 3).Value = "18:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = 9
$ws.Cells.Item(67, 6).Value = "Occupied"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 1).ClearFormats()
$ws.Cells.Item(68, 2).Value = "18:32:28"
$ws.Cells.Item(68, 3).Value = "18:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = 2
$ws.Cells.Item(68, 6).Value = "Occupied"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 1).ClearFormats()
$ws.Cells.Item(69, 2).Value = "18:32:31"
$ws.Cells.Item(69, 3).Value = "18:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = 1
$ws.Cells.Item(69, 6).Value = "Occupied"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 1).ClearFormats()
$ws.Cells.Item(70, 2).Value = "18:32:35"
$ws.Cells.Item(70, 3).Value = "18:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = 2
$ws.Cells.Item(70, 6).Value = "Occupied"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 1).ClearFormats()
$ws.Cells.Item(71, 2).Value = "18:32:38"
$ws.Cells.Item(71, 3).Value = "18:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = 1
$ws.Cells.Item(71, 6).Value = "Occupied"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(72, 2).Value = "18:32:41"
$ws.Cells.Item(72, 3).Value = "18:00"
$ws.Cells.Item(72, 4).Value = "Bedroom"
$ws.Cells.Item(72, 5).Value = 2
$ws.Cells.Item(72, 6).Value = "Occupied"
$ws.Cells.Item(73, 1).Value = "'2026-01-28"
$ws.Cells.Item(73, 1).ClearFormats()
$ws.Cells.Item(73, 2).Value = "18:32:44"
$ws.Cells.Item(73, 3).Value = "18:00"
$ws.Cells.Item(73, 4).Value = "Bedroom"
$ws.Cells.Item(73, 5).Value = 1
$ws.Cells.Item(73, 6).Value = "Occupied"

$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Cells.Item(58, 1).Value = "'2026-01-28"
$ws.Cells.Item(58, 1).ClearFormats()
$ws.Cells.Item(58, 2).Value = "18:31:54"
$ws.Cells.Item(58, 3).Value = "18:00"
$ws.Cells.Item(58, 4).Value = "Bedroom"
$ws.Cells.Item(58, 5).Value = "In Bed"
$ws.Cells.Item(58, 6).Value = "Occupied"
$ws.Cells.Item(59, 1).Value = "'2026-01-28"
$ws.Cells.Item(59, 1).ClearFormats()
$ws.Cells.Item(59, 2).Value = "18:31:57"
$ws.Cells.Item(59, 3).Value = "18:00"
$ws.Cells.Item(59, 4).Value = "Bedroom"
$ws.Cells.Item(59, 5).Value = "In Bed"
$ws.Cells.Item(59, 6).Value = "Occupied"
$ws.Cells.Item(60, 1).Value = "'2026-01-28"
$ws.Cells.Item(60, 1).ClearFormats()
$ws.Cells.Item(60, 2).Value = "18:32:00"
$ws.Cells.Item(60, 3).Value = "18:00"
$ws.Cells.Item(60, 4).Value = "Bedroom"
$ws.Cells.Item(60, 5).Value = "In Bed"
$ws.Cells.Item(60, 6).Value = "Occupied"
$ws.Cells.Item(61, 1).Value = "'2026-01-28"
$ws.Cells.Item(61, 1).ClearFormats()
$ws.Cells.Item(61, 2).Value = "18:32:04"
$ws.Cells.Item(61, 3).Value = "18:00"
$ws.Cells.Item(61, 4).Value = "Bedroom"
$ws.Cells.Item(61, 5).Value = "In Bed"
$ws.Cells.Item(61, 6).Value = "Occupied"
$ws.Cells.Item(62, 1).Value = "'2026-01-28"
$ws.Cells.Item(62, 1).ClearFormats()
$ws.Cells.Item(62, 2).Value = "18:32:08"
$ws.Cells.Item(62, 3).Value = "18:00"
$ws.Cells.Item(62, 4).Value = "Bedroom"
$ws.Cells.Item(62, 5).Value = "In Bed"
$ws.Cells.Item(62, 6).Value = "Occupied"
$ws.Cells.Item(63, 1).Value = "'2026-01-28"
$ws.Cells.Item(63, 1).ClearFormats()
$ws.Cells.Item(63, 2).Value = "18:32:11"
$ws.Cells.Item(63, 3).Value = "18:00"
$ws.Cells.Item(63, 4).Value = "Bedroom"
$ws.Cells.Item(63, 5).Value = "In Bed"
$ws.Cells.Item(63, 6).Value = "Occupied"
$ws.Cells.Item(64, 1).Value = "'2026-01-28"
$ws.Cells.Item(64, 1).ClearFormats()
$ws.Cells.Item(64, 2).Value = "18:32:14"
$ws.Cells.Item(64, 3).Value = "18:00"
$ws.Cells.Item(64, 4).Value = "Bedroom"
$ws.Cells.Item(64, 5).Value = "In Bed"
$ws.Cells.Item(64, 6).Value = "Occupied"
$ws.Cells.Item(65, 1).Value = "'2026-01-28"
$ws.Cells.Item(65, 1).ClearFormats()
$ws.Cells.Item(65, 2).Value = "18:32:17"
$ws.Cells.Item(65, 3).Value = "18:00"
$ws.Cells.Item(65, 4).Value = "Bedroom"
$ws.Cells.Item(65, 5).Value = "In Bed"
$ws.Cells.Item(65, 6).Value = "Occupied"
$ws.Cells.Item(66, 1).Value = "'2026-01-28"
$ws.Cells.Item(66, 1).ClearFormats()
$ws.Cells.Item(66, 2).Value = "18:32:20"
$ws.Cells.Item(66, 3).Value = "18:00"
$ws.Cells.Item(66, 4).Value = "Bedroom"
$ws.Cells.Item(66, 5).Value = "In Bed"
$ws.Cells.Item(66, 6).Value = "Occupied"
$ws.Cells.Item(67, 1).Value = "'2026-01-28"
$ws.Cells.Item(67, 1).ClearFormats()
$ws.Cells.Item(67, 2).Value = "18:32:23"
$ws.Cells.Item(67, 3).Value = "18:00"
$ws.Cells.Item(67, 4).Value = "Bedroom"
$ws.Cells.Item(67, 5).Value = "In Bed"
$ws.Cells.Item(67, 6).Value = "Occupied"
$ws.Cells.Item(68, 1).Value = "'2026-01-28"
$ws.Cells.Item(68, 1).ClearFormats()
$ws.Cells.Item(68, 2).Value = "18:32:26"
$ws.Cells.Item(68, 3).Value = "18:00"
$ws.Cells.Item(68, 4).Value = "Bedroom"
$ws.Cells.Item(68, 5).Value = "In Bed"
$ws.Cells.Item(68, 6).Value = "Occupied"
$ws.Cells.Item(69, 1).Value = "'2026-01-28"
$ws.Cells.Item(69, 1).ClearFormats()
$ws.Cells.Item(69, 2).Value = "18:32:29"
$ws.Cells.Item(69, 3).Value = "18:00"
$ws.Cells.Item(69, 4).Value = "Bedroom"
$ws.Cells.Item(69, 5).Value = "In Bed"
$ws.Cells.Item(69, 6).Value = "Occupied"
$ws.Cells.Item(70, 1).Value = "'2026-01-28"
$ws.Cells.Item(70, 1).ClearFormats()
$ws.Cells.Item(70, 2).Value = "18:32:33"
$ws.Cells.Item(70, 3).Value = "18:00"
$ws.Cells.Item(70, 4).Value = "Bedroom"
$ws.Cells.Item(70, 5).Value = "In Bed"
$ws.Cells.Item(70, 6).Value = "Occupied"
$ws.Cells.Item(71, 1).Value = "'2026-01-28"
$ws.Cells.Item(71, 1).ClearFormats()
$ws.Cells.Item(71, 2).Value = "18:32:36"
$ws.Cells.Item(71, 3).Value = "18:00"
$ws.Cells.Item(71, 4).Value = "Bedroom"
$ws.Cells.Item(71, 5).Value = "In Bed"
$ws.Cells.Item(71, 6).Value = "Occupied"
$ws.Cells.Item(72, 1).Value = "'2026-01-28"
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(72, 2).Value = "18:32:39"
$ws.Cells.Item(72, 3).Value = "18:00"
$ws.Cells.Item(72, 4).Value = "Bedroom"
$ws.Cells.Item(72, 5).Value = "In Bed"
$ws.Cells.Item(72, 6).Value = "Occupied"
$ws.Cells.Item(73, 1).Value = "'2026-01-28"
$ws.Cells.Item(73, 1).ClearFormats()
$ws.Cells.Item(73, 2).Value = "18:32:42"
$ws.Cells.Item(73, 3).Value = "18:00"
$ws.Cells.Item(73, 4).Value = "Bedroom"
$ws.Cells.Item(73, 5).Value = "In Bed"
$ws.Cells.Item(73, 6).Value = "Occupied"

$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(21, 1).Value = "'2026-01-28"
$ws.Cells.Item(21, 1).ClearFormats()
$ws.Cells.Item(21, 2).Value = "18:32:03"
$ws.Cells.Item(21, 3).Value = "18:00"
$ws.Cells.Item(21, 4).Value = "Living Room Main Door"
$ws.Cells.Item(21, 5).Value = "ENTER"
$ws.Cells.Item(21, 6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(22, 1).Value = "'2026-01-28"
$ws.Cells.Item(22, 1).ClearFormats()
$ws.Cells.Item(22, 2).Value = "18:32:32"
$ws.Cells.Item(22, 3).Value = "18:00"
$ws.Cells.Item(22, 4).Value = "Living Room Main Door"
$ws.Cells.Item(22, 5).Value = "EXIT"
$ws.Cells.Item(22, 6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(23, 1).Value = "'2026-01-28"
$ws.Cells.Item(23, 1).ClearFormats()
$ws.Cells.Item(23, 2).Value = "18:32:45"
$ws.Cells.Item(23, 3).Value = "18:00"
$ws.Cells.Item(23, 4).Value = "Living Room Main Door"
$ws.Cells.Item(23, 5).Value = "ENTER"
$ws.Cells.Item(23, 6).Value = "User ENTERED Living Room Main Door"

$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(13, 1).Value = "'2026-01-28"
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(13, 2).Value = "18:32:07"
$ws.Cells.Item(13, 3).Value = "18:00"
$ws.Cells.Item(13, 4).Value = "Living Room Main Door"
$ws.Cells.Item(13, 5).Value = "Image Captured"
$ws.Cells.Item(13, 6).Value = "Active"
$ws.Cells.Item(14, 1).Value = "'2026-01-28"
$ws.Cells.Item(14, 1).ClearFormats()
$ws.Cells.Item(14, 2).Value = "18:32:46"
$ws.Cells.Item(14, 3).Value = "18:00"
$ws.Cells.Item(14, 4).Value = "Living Room Main Door"
$ws.Cells.Item(14, 5).Value = "Image Captured"
$ws.Cells.Item(14, 6).Value = "Active"

